$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Date as Text:" example value (C10)
$ws.Range("C10").Value = "9/2/2010 12:00:00 AM"

# Update the "DateTime as Text:" example value (C11)
$ws.Range("C11").Value = "9/2/2010 1:45:22 PM"

# Update the "DateTime to Text:" example value (C21)
$ws.Range("C21").Value = "9/2/10 1:45"

# Update the "Formatted Number to Text:" example value (C33).
# This cell is styled with a numeric format, so temporarily switch it to
# text format while assigning the new text so Excel doesn't re-parse the
# string back into a number; then restore the original format.
$origFormat = $ws.Range("C33").NumberFormat
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = "12,345.68"
$ws.Range("C33").NumberFormat = $origFormat

# Widen column C slightly
$ws.Columns.Item(3).ColumnWidth = 21.7
